$wb = $excel.ActiveWorkbook

# --- "TODO" sheet: add a "DONE" marker column next to the Expiry Date task rows ---
$todo = $wb.Worksheets.Item("TODO")
$todo.Range("C3").Value = "DONE"
$todo.Range("C4").Value = "DONE"
$todo.Range("C5").Value = "DONE"

$todo.Activate()
$todo.Range("C5").Select()

# --- "11 JAN 2017" sheet: move the saved selection (no data change) ---
$bugs = $wb.Worksheets.Item("11 JAN 2017")
$bugs.Activate()
$bugs.Range("B37").Select()

# leave TODO as the active sheet, matching the saved workbook view (activeTab stays on TODO)
$todo.Activate()
